$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A32").Value = 31
$ws.Range("B32").Value = "BRAND NEW WOWAOWOOWAOAO"
$ws.Range("C32").Value = "pelae"
$ws.Range("D32").Value = "work"
$ws.Range("E32").Value = "ples"
$ws.Range("F32").Value = "yes"

$ws.Range("F32").Select()
